# IRYO-vaccination_data.xlsx — add the 2021-05-19 daily row.
#
# The source sheet keeps a reverse-chronological (newest-first) log of daily
# vaccination counts starting at row 5 (row 4 is the running total). A new
# day's figures are published by inserting a fresh row right under the
# header/total row and pushing all the existing daily rows down by one,
# then updating the grand-total row and the "as of" date caption.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 5 — everything currently at row 5.. moves
#    down to row 6.. (formulas/shared-formula ranges, styles and the
#    calcChain all shift automatically, exactly like a manual
#    right-click > Insert in Excel).
$ws.Rows("5:5").Insert()

# 2. The new row 5 has no formatting yet; pull it from row 6 (which now
#    holds what used to be row 5, i.e. an identically-formatted daily row)
#    so the new row matches the rest of the table (date format, borders,
#    centred weekday text, number formatting, ...).
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Fill in the new day's data: 2021-05-19 (serial 44335), a Wednesday.
$ws.Range("A5").Value = 44335
$ws.Range("B5").Value = "(水)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 88163
$ws.Range("E5").Value = 179703

# 4. Update the grand-total row (row 4) components; the C4 total is a
#    formula (=SUM(D4:E4)) and recalculates on its own.
$ws.Range("D4").Value = 3784071
$ws.Range("E4").Value = 2214762

# 5. Bump the "as of" caption from 5/18 to 5/19.
$ws.Range("E2").Value = "（5月19日時点）"
